$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Summary")

# These cells hold numeric-looking figures that are stored as text (shared
# strings) in the original workbook, with the default "Normal" / General
# style (no explicit number format). Force a Text format just long enough
# to assign the new text value (so COM doesn't silently convert the string
# into a real number), then restore the original "Normal" style so the
# cell formatting is left exactly as it was.

# Enterprises density (per 1000 people) - SMEs column: 3.8 -> 3.81
$ws.Range("C11").NumberFormat = "@"
$ws.Range("C11").Value = "3.81"
$ws.Range("C11").Style = "Normal"

# Enterprises (% of total): Micro 44.5 -> 44.48, SMEs 53 -> 52.99, MSMEs 97.5 -> 97.46
$ws.Range("B12").NumberFormat = "@"
$ws.Range("B12").Value = "44.48"
$ws.Range("B12").Style = "Normal"

$ws.Range("C12").NumberFormat = "@"
$ws.Range("C12").Value = "52.99"
$ws.Range("C12").Style = "Normal"

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "97.46"
$ws.Range("D12").Style = "Normal"
